$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the numeric-looking Price/Volume columns as Text first so COM
# keeps the assigned strings verbatim (matches the inlineStr source data)
# instead of auto-converting "314.51" / "0.82%" into numbers.
$numericRange = $ws.Range("D2:E51")
$numericRange.NumberFormat = "@"

$ws.Range("D2").Value = "314.51"
$ws.Range("E2").Value = "0.82%"

$ws.Range("D3").Value = "48.51"
$ws.Range("E3").Value = "8.94%"

$ws.Range("D4").Value = "5.289"
$ws.Range("E4").Value = "3.79%"

$ws.Range("D5").Value = "0.07887"
$ws.Range("E5").Value = "-1.60%"

$ws.Range("D6").Value = "4.588"
$ws.Range("E6").Value = "2.27%"

$ws.Range("D7").Value = "1.321"
$ws.Range("E7").Value = "22.44%"

$ws.Range("D8").Value = "1.613"
$ws.Range("E8").Value = "-2.25%"

$ws.Range("D9").Value = "0.1236"
$ws.Range("E9").Value = "-4.07%"

$ws.Range("D10").Value = "0.1949"
$ws.Range("E10").Value = "2.93%"

$ws.Range("D11").Value = "0.09560"
$ws.Range("E11").Value = "3.56%"

$ws.Range("D12").Value = "0.04546"
$ws.Range("E12").Value = "8.28%"

$ws.Range("D13").Value = "0.1049"
$ws.Range("E13").Value = "1.11%"

$ws.Range("D14").Value = "0.001307"
$ws.Range("E14").Value = "0.04%"

$ws.Range("D15").Value = "0.04216"
$ws.Range("E15").Value = "0.58%"

$ws.Range("D16").Value = "0.005860"
$ws.Range("E16").Value = "0.21%"

$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.335"
$ws.Range("E17").Value = "-1.10%"

$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "2.468"
$ws.Range("E18").Value = "2.77%"

$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "0.3464"
$ws.Range("E19").Value = "3.05%"

$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").Value = "8.060"
$ws.Range("E20").Value = "0.81%"

$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "0.1394"
$ws.Range("E21").Value = "1.19%"

$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "0.3074"
$ws.Range("E22").Value = "-1.78%"

$ws.Range("B23").Value = "BitKan"
$ws.Range("C23").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D23").Value = "0.001298"
$ws.Range("E23").Value = "2.18%"

$ws.Range("B24").Value = "HotbitToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D24").Value = "0.004181"
$ws.Range("E24").Value = "-9.15%"

$ws.Range("D25").Value = "0.0001358"
$ws.Range("E25").Value = "1.65%"

$ws.Range("E26").Value = "-95.20%"

$ws.Range("D38").Value = "0.02647"
$ws.Range("E38").Value = "-0.35%"

$ws.Range("D39").Value = "0.05802"
$ws.Range("E39").Value = "7.07%"

$ws.Range("D40").Value = "0.01082"
$ws.Range("E40").Value = "92.97%"

$ws.Range("D41").Value = "0.008001"
$ws.Range("E41").Value = "3.68%"

$ws.Range("D42").Value = "0.1442"
$ws.Range("E42").Value = "2.12%"

$ws.Range("D43").Value = "0.008162"
$ws.Range("E43").Value = "11.77%"

$ws.Range("D44").Value = "0.008669"
$ws.Range("E44").Value = "3.37%"

$ws.Range("D45").Value = "0.3146"
$ws.Range("E45").Value = "0.76%"

$ws.Range("D46").Value = "0.00007015"
$ws.Range("E46").Value = "4.66%"

$ws.Range("E47").Value = "1.62%"

$ws.Range("D48").Value = "0.05599"
$ws.Range("E48").Value = "5.57%"

$ws.Range("D49").Value = "0.004023"
$ws.Range("E49").Value = "1.65%"

$ws.Range("E50").Value = "1.62%"

$ws.Range("D51").Value = "0.0002011"
$ws.Range("E51").Value = "1.62%"

# Restore the default style on the touched numeric range so only the
# values change (no leftover number-format/style delta).
$numericRange.Style = "Normal"
